# Applies the Pflichtenheft.docx edit:
#  1. Drop the stray "_GoBack" bookmark that sat right after "IP-Telefonie 1.1"
#     on the title page.
#  2. Reword the MUSS-Kriterien bullet about SRTP, and re-create the
#     "_GoBack" bookmark in its new spot, between "abhörbar" and " sein".

$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark (title page, after "IP-Telefonie 1.1").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Swap the old SRTP sentence for the new wording.
$d.Content.Find.Execute(
    "Es sollte die interne Kommunikation per SRTP abgesichert werden.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Die Gespräche dürfen nicht abhörbar sein", 2)

# 3. Re-create "_GoBack" right between "abhörbar" and " sein", matching the
#    run split in the target document.
$r = $d.Content
$r.Find.Execute("abhörbar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
